# Fix "Ellenberg N" / "Ellenberg L" labels that were typed with a
# non-breaking space (U+00A0) between the words instead of a normal
# space (U+0020). Also restore the active sheet / selection to the
# "relasjoner" sheet, matching where the fix was made.

$wb = $excel.ActiveWorkbook

$ellenbergN = "Ellenberg" + [char]0x20 + "N"
$ellenbergL = "Ellenberg" + [char]0x20 + "L"

$ws1 = $wb.Worksheets.Item("tilstandsindikatorer")
$ws5 = $wb.Worksheets.Item("relasjoner")

# tilstandsindikatorer!B12 -> "Ellenberg L"
$ws1.Range("B12").Value2 = $ellenbergL

# relasjoner!C11, B32, C51 -> "Ellenberg N"
$ws5.Range("C11").Value2 = $ellenbergN
$ws5.Range("B32").Value2 = $ellenbergN
$ws5.Range("C51").Value2 = $ellenbergN

# relasjoner!C12, B33, C52, C68 -> "Ellenberg L"
$ws5.Range("C12").Value2 = $ellenbergL
$ws5.Range("B33").Value2 = $ellenbergL
$ws5.Range("C52").Value2 = $ellenbergL
$ws5.Range("C68").Value2 = $ellenbergL

# The author last worked on the "relasjoner" sheet, so it is the active
# tab/sheet, with a specific cell selected, when the file was saved.
$ws5.Activate()
$ws5.Range("C21").Select()
